$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- ws1 new rows ---
# row 3
$ws1.Cells.Item(3,1).Value = "20230225RA2502"
$ws1.Cells.Item(3,2).Value = 5000
$ws1.Cells.Item(3,3).Value = 0.05
$ws1.Cells.Item(3,4).Value = 12
$ws1.Cells.Item(3,5).NumberFormat = "@"
$ws1.Cells.Item(3,5).Value = "2023-03-09"
$ws1.Cells.Item(3,5).Style = "Normal"
$ws1.Cells.Item(3,6).NumberFormat = "@"
$ws1.Cells.Item(3,6).Value = "2023-03-09"
$ws1.Cells.Item(3,6).Style = "Normal"
$ws1.Cells.Item(3,7).NumberFormat = "@"
$ws1.Cells.Item(3,7).Value = "2024-03-09"
$ws1.Cells.Item(3,7).Style = "Normal"
$ws1.Cells.Item(3,8).Value = 437.5
$ws1.Cells.Item(3,9).Value = 5250
$ws1.Cells.Item(3,10).Value = "Rahul"
$ws1.Cells.Item(3,11).NumberFormat = "@"
$ws1.Cells.Item(3,11).Value = "0.03"
$ws1.Cells.Item(3,11).Style = "Normal"
$ws1.Cells.Item(3,12).Value = 1.575
$ws1.Cells.Item(3,13).Value = 12
$ws1.Cells.Item(3,15).Value = "pending"
$ws1.Cells.Item(3,16).Value = "pending"
$ws1.Cells.Item(3,17).Value = "pending"
$ws1.Cells.Item(3,18).Value = "pending"
$ws1.Cells.Item(3,19).Value = "pending"
$ws1.Cells.Item(3,20).Value = "pending"
$ws1.Cells.Item(3,21).Value = "pending"
$ws1.Cells.Item(3,22).Value = "pending"
$ws1.Cells.Item(3,23).Value = "pending"
$ws1.Cells.Item(3,24).Value = "pending"
$ws1.Cells.Item(3,25).Value = "pending"
$ws1.Cells.Item(3,26).Value = "pending"

# row 4
$ws1.Cells.Item(4,1).Value = "20230317DE1703"
$ws1.Cells.Item(4,2).Value = 120000
$ws1.Cells.Item(4,3).Value = 0.5
$ws1.Cells.Item(4,4).Value = 12
$ws1.Cells.Item(4,5).NumberFormat = "yyyy\-mm\-dd"
$ws1.Cells.Item(4,5).Value = 45002
$ws1.Cells.Item(4,6).NumberFormat = "yyyy\-mm\-dd"
$ws1.Cells.Item(4,6).Value = 45002
$ws1.Cells.Item(4,7).NumberFormat = "yyyy\-mm\-dd"
$ws1.Cells.Item(4,7).Value = 45368
$ws1.Cells.Item(4,8).Value = 10050
$ws1.Cells.Item(4,9).Value = 120600
$ws1.Cells.Item(4,10).Value = "Rahul"
$ws1.Cells.Item(4,11).Value = 0.1
$ws1.Cells.Item(4,12).Value = 120.6
$ws1.Cells.Item(4,13).Value = 12
$ws1.Cells.Item(4,15).Value = "pending"
$ws1.Cells.Item(4,16).Value = "pending"
$ws1.Cells.Item(4,17).Value = "pending"
$ws1.Cells.Item(4,18).Value = "pending"
$ws1.Cells.Item(4,19).Value = "pending"
$ws1.Cells.Item(4,20).Value = "pending"
$ws1.Cells.Item(4,21).Value = "pending"
$ws1.Cells.Item(4,22).Value = "pending"
$ws1.Cells.Item(4,23).Value = "pending"
$ws1.Cells.Item(4,24).Value = "pending"
$ws1.Cells.Item(4,25).Value = "pending"
$ws1.Cells.Item(4,26).Value = "pending"

# --- ws2 new rows ---
# row 3
$ws2.Cells.Item(3,1).Value = "20230225RA2502"
$ws2.Cells.Item(3,2).NumberFormat = "yyyy\-mm\-dd"
$ws2.Cells.Item(3,2).Value = 44994
$ws2.Cells.Item(3,3).NumberFormat = "yyyy\-mm\-dd"
$ws2.Cells.Item(3,3).Value = 45025
$ws2.Cells.Item(3,4).NumberFormat = "yyyy\-mm\-dd"
$ws2.Cells.Item(3,4).Value = 45055
$ws2.Cells.Item(3,5).NumberFormat = "yyyy\-mm\-dd"
$ws2.Cells.Item(3,5).Value = 45086
$ws2.Cells.Item(3,6).NumberFormat = "yyyy\-mm\-dd"
$ws2.Cells.Item(3,6).Value = 45116
$ws2.Cells.Item(3,7).NumberFormat = "yyyy\-mm\-dd"
$ws2.Cells.Item(3,7).Value = 45147
$ws2.Cells.Item(3,8).NumberFormat = "yyyy\-mm\-dd"
$ws2.Cells.Item(3,8).Value = 45178
$ws2.Cells.Item(3,9).NumberFormat = "yyyy\-mm\-dd"
$ws2.Cells.Item(3,9).Value = 45208
$ws2.Cells.Item(3,10).NumberFormat = "yyyy\-mm\-dd"
$ws2.Cells.Item(3,10).Value = 45239
$ws2.Cells.Item(3,11).NumberFormat = "yyyy\-mm\-dd"
$ws2.Cells.Item(3,11).Value = 45269
$ws2.Cells.Item(3,12).NumberFormat = "yyyy\-mm\-dd"
$ws2.Cells.Item(3,12).Value = 45300
$ws2.Cells.Item(3,13).NumberFormat = "yyyy\-mm\-dd"
$ws2.Cells.Item(3,13).Value = 45331

# row 4
$ws2.Cells.Item(4,1).Value = "20230317DE1703"
$ws2.Cells.Item(4,2).NumberFormat = "yyyy\-mm\-dd"
$ws2.Cells.Item(4,2).Value = 44996
$ws2.Cells.Item(4,3).NumberFormat = "yyyy\-mm\-dd"
$ws2.Cells.Item(4,3).Value = 45033
$ws2.Cells.Item(4,4).NumberFormat = "yyyy\-mm\-dd"
$ws2.Cells.Item(4,4).Value = 45063
$ws2.Cells.Item(4,5).NumberFormat = "yyyy\-mm\-dd"
$ws2.Cells.Item(4,5).Value = 45094
$ws2.Cells.Item(4,6).NumberFormat = "yyyy\-mm\-dd"
$ws2.Cells.Item(4,6).Value = 45124
$ws2.Cells.Item(4,7).NumberFormat = "yyyy\-mm\-dd"
$ws2.Cells.Item(4,7).Value = 45155
$ws2.Cells.Item(4,8).NumberFormat = "yyyy\-mm\-dd"
$ws2.Cells.Item(4,8).Value = 45186
$ws2.Cells.Item(4,9).NumberFormat = "yyyy\-mm\-dd"
$ws2.Cells.Item(4,9).Value = 45216
$ws2.Cells.Item(4,10).NumberFormat = "yyyy\-mm\-dd"
$ws2.Cells.Item(4,10).Value = 45247
$ws2.Cells.Item(4,11).NumberFormat = "yyyy\-mm\-dd"
$ws2.Cells.Item(4,11).Value = 45277
$ws2.Cells.Item(4,12).NumberFormat = "yyyy\-mm\-dd"
$ws2.Cells.Item(4,12).Value = 45308
$ws2.Cells.Item(4,13).NumberFormat = "yyyy\-mm\-dd"
$ws2.Cells.Item(4,13).Value = 45339
